$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (16) of data to the sheet, mirroring the existing rows.
$row = 16

$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item($row, 1).Value = 42622.887256944443

$ws.Cells.Item($row, 2).Value = -22
$ws.Cells.Item($row, 3).Value = 65
$ws.Cells.Item($row, 4).Value = 33
$ws.Cells.Item($row, 5).Value = 65
$ws.Cells.Item($row, 6).Value = 18
$ws.Cells.Item($row, 7).Value = 15072
$ws.Cells.Item($row, 8).Value = 10502
$ws.Cells.Item($row, 9).Value = 1632
$ws.Cells.Item($row, 10).Value = 206
$ws.Cells.Item($row, 11).Value = 104
$ws.Cells.Item($row, 12).Value = 18
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Noun"
